$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values look numeric, to preserve them as text
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply the updated values
$ws.Range('D2').Value = '57.052.47'
$ws.Range('E2').Value = '  +0.77%  '
$ws.Range('D3').Value = '2.428.59'
$ws.Range('E3').Value = '  -2.45%  '
$ws.Range('E4').Value = '  +0.29%  '
$ws.Range('D5').Value = '486.24'
$ws.Range('E5').Value = '  -1.11%  '
$ws.Range('D6').Value = '154.54'
$ws.Range('E6').Value = '  +2.19%  '
$ws.Range('D7').Value = '0.616'
$ws.Range('E7').Value = '  +19.58%  '
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  +0.22%  '
$ws.Range('D9').Value = '2.439.72'
$ws.Range('E9').Value = '  -2.43%  '
$ws.Range('D10').Value = '0.0992'
$ws.Range('E10').Value = '  +0.96%  '
$ws.Range('E11').Value = '  -2.08%  '
$ws.Range('E12').Value = '  -0.46%  '
$ws.Range('E13').Value = '  +0.79%  '
$ws.Range('D14').Value = '2.852.71'
$ws.Range('E14').Value = '  -2.10%  '
$ws.Range('D15').Value = '57.170.72'
$ws.Range('E15').Value = '  +0.90%  '
$ws.Range('D16').Value = '20.72'
$ws.Range('E16').Value = '  -2.30%  '
$ws.Range('D17').Value = '0.0000133'
$ws.Range('E17').Value = '  -2.45%  '
$ws.Range('D18').Value = '2.437.16'
$ws.Range('E18').Value = '  -2.56%  '
$ws.Range('D19').Value = '4.74'
$ws.Range('E19').Value = '  +4.53%  '
$ws.Range('D20').Value = '326.51'
$ws.Range('E20').Value = '  +1.29%  '
$ws.Range('D21').Value = '9.96'
$ws.Range('E21').Value = '  -2.95%  '
$ws.Range('E22').Value = '  +0.22%  '
$ws.Range('D23').Value = '5.89'
$ws.Range('E23').Value = '  -0.20%  '
$ws.Range('D24').Value = '58.19'
$ws.Range('E24').Value = '  -1.02%  '
$ws.Range('D25').Value = '0.409'
$ws.Range('E25').Value = '  -0.36%  '
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  +0.18%  '
$ws.Range('E27').Value = '  -3.02%  '
$ws.Range('D28').Value = '2.544.11'
$ws.Range('E28').Value = '  -2.13%  '
$ws.Range('D29').Value = '7.22'
$ws.Range('E29').Value = '  -5.71%  '
$ws.Range('D30').Value = '0.0₃0779'
$ws.Range('E30').Value = '  -2.57%  '
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  +0.23%  '
$ws.Range('D32').Value = '18.71'
$ws.Range('E32').Value = '  +1.73%  '
$ws.Range('D33').Value = '148.90'
$ws.Range('E33').Value = '  -1.41%  '
$ws.Range('D34').Value = '1.52'
$ws.Range('E34').Value = '  +0.07%  '
$ws.Range('D35').Value = '5.30'
$ws.Range('E35').Value = '  +1.86%  '
$ws.Range('E36').Value = '  -1.42%  '
$ws.Range('D37').Value = '3.69'
$ws.Range('E37').Value = '  -2.08%  '
$ws.Range('D38').Value = '0.849'
$ws.Range('E38').Value = '  -2.89%  '
$ws.Range('D39').Value = '0.102'
$ws.Range('E39').Value = '  +10.46%  '
$ws.Range('D40').Value = '34.22'
$ws.Range('E40').Value = '  +0.81%  '
$ws.Range('E41').Value = '  -1.14%  '
$ws.Range('D42').Value = '3.51'
$ws.Range('E42').Value = '  +0.02%  '
$ws.Range('E43').Value = '  +0.20%  '
$ws.Range('D44').Value = '0.595'
$ws.Range('E44').Value = '  -2.84%  '
$ws.Range('D45').Value = '265.61'
$ws.Range('E45').Value = '  +0.47%  '
$ws.Range('D46').Value = '0.0531'
$ws.Range('E46').Value = '  -4.83%  '
$ws.Range('D47').Value = '10.20'
$ws.Range('E47').Value = '  -0.25%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').Value = '0.0227'
$ws.Range('E48').Value = '  -1.28%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').Value = '4.66'
$ws.Range('E49').Value = '  -4.05%  '
$ws.Range('D50').Value = '17.49'
$ws.Range('E50').Value = '  -1.56%  '
$ws.Range('D51').Value = '1.858.46'
$ws.Range('E51').Value = '  -2.42%  '
